$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (Player, Position, Team) replacing rows 2-19
$data = @(
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Ziaire Williams", "SG,SF", "Brooklyn Nets"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Max Christie", "SG,SF", "Dallas Mavericks"),
    @("Jimmy Butler III", "SF,PF", "Golden State Warriors"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Brandon Ingram", "SG,SF,PF", "Toronto Raptors"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
